$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before column J ("clip_ratio" | "vf_lr" | "target_kl")
$ws.Columns.Item(10).Insert()

# 2) New header for the inserted column
$ws.Cells.Item(1, 10).Value = "vf_lr"

# 3) Fill the new vf_lr column for the existing trials (rows 2-19) with 0.001
$ws.Range("J2:J19").Value = 0.001

# 4) Format the new column like the other learning-rate style columns:
#    scientific notation, centered
$ws.Range("J2:J19").NumberFormat = "0.00E+00"
$ws.Range("J2:J19").HorizontalAlignment = -4108

# 5) Add the two new trial rows (15 second run, and new Trial 16)
$ws.Range("A20").Value = 15
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = "Hallway (setups)"
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0.5
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 128
$ws.Range("H20").Value = 2
$ws.Range("I20").Value = 0.1
$ws.Range("J20").Value = 0.001
$ws.Range("K20").Value = 0.015
$ws.Range("L20").Value = 0.99
$ws.Range("M20").Value = 0.95
$ws.Range("N20").Value = 4000
$ws.Range("O20").Value = 800
$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = "Yes"
$ws.Range("R20").Value = "Bad"

$ws.Range("A21").Value = 16
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "Hallway (setups)"
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 0.5
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 128
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = 0.1
$ws.Range("J21").Value = 0.0005
$ws.Range("K21").Value = 0.015
$ws.Range("L21").Value = 0.99
$ws.Range("M21").Value = 0.95
$ws.Range("N21").Value = 4000
$ws.Range("O21").Value = 1000
$ws.Range("P21").Value = 2
$ws.Range("Q21").Value = "Yes"
$ws.Range("R21").Value = "Bad at 200 and 1000"

$ws.Range("J20:J21").NumberFormat = "0.00E+00"
$ws.Range("J20:J21").HorizontalAlignment = -4108

# 6) Update the view: scroll right a bit and move the selection like the saved file
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("S13").Select()
